# Auto-generated script to apply numeric value updates to the Sheets workbook
# (values refreshed by the scheduled market-data runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5385.643
$ws.Range("I74").Value = 5189.9
$ws.Range("J74").Value = 5875
$ws.Range("K74").Value = 5189.9
$ws.Range("L74").Value = 5875
$ws.Range("M74").Value = -4253.9
$ws.Range("N74").Value = -7747
$ws.Range("H76").Value = 3524.8333
$ws.Range("I76").Value = 3269.9
$ws.Range("J76").Value = 4799.5
$ws.Range("K76").Value = 3269.9
$ws.Range("L76").Value = 4799.5
$ws.Range("M76").Value = -2954.9
$ws.Range("N76").Value = -5429.5
$ws.Range("H77").Value = 5385.643
$ws.Range("I77").Value = 5189.9
$ws.Range("J77").Value = 5875
$ws.Range("K77").Value = 25949.5
$ws.Range("L77").Value = 29375
$ws.Range("M77").Value = -21269.5
$ws.Range("N77").Value = -38735
$ws.Range("H79").Value = 3524.8333
$ws.Range("I79").Value = 3269.9
$ws.Range("J79").Value = 4799.5
$ws.Range("K79").Value = 3269.9
$ws.Range("L79").Value = 4799.5
$ws.Range("M79").Value = -2177.9
$ws.Range("N79").Value = -6983.5
$ws.Range("H92").Value = 432.61905
$ws.Range("I92").Value = 432.61905
$ws.Range("K92").Value = 432.61905
$ws.Range("M92").Value = 815.38095
$ws.Range("H113").Value = 102450.5
$ws.Range("I113").Value = 252251.25
$ws.Range("J113").Value = 2583.3333
$ws.Range("K113").Value = 252251.25
$ws.Range("L113").Value = 2583.3333
$ws.Range("M113").Value = -248997.25
$ws.Range("N113").Value = -9091.3333
$ws.Range("H137").Value = 981.4857
$ws.Range("I137").Value = 869.7037
$ws.Range("J137").Value = 1358.75
$ws.Range("K137").Value = 2609.1111
$ws.Range("L137").Value = 4076.25
$ws.Range("M137").Value = -59.11110000000008
$ws.Range("N137").Value = -9176.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2694.889
$ws.Range("I63").Value = 2237.182
$ws.Range("J63").Value = 3414.1428
$ws.Range("K63").Value = 2237.182
$ws.Range("L63").Value = 3414.1428
$ws.Range("M63").Value = -1551.182
$ws.Range("N63").Value = -4786.1428
$ws.Range("H66").Value = 2694.889
$ws.Range("I66").Value = 2237.182
$ws.Range("J66").Value = 3414.1428
$ws.Range("K66").Value = 11185.91
$ws.Range("L66").Value = 17070.714
$ws.Range("M66").Value = -7753.91
$ws.Range("N66").Value = -23934.714
$ws.Range("H102").Value = 68702
$ws.Range("I102").Value = 92961.82000000001
$ws.Range("J102").Value = 1987.5
$ws.Range("K102").Value = 92961.82000000001
$ws.Range("L102").Value = 1987.5
$ws.Range("M102").Value = -91339.82000000001
$ws.Range("N102").Value = -5231.5
$ws.Range("H112").Value = 17166.666
$ws.Range("J112").Value = 17166.666
$ws.Range("L112").Value = 17166.666
$ws.Range("N112").Value = -20120.666
$ws.Range("H122").Value = 2260.6843
$ws.Range("I122").Value = 2710.3635
$ws.Range("J122").Value = 1642.375
$ws.Range("K122").Value = 8131.0905
$ws.Range("L122").Value = 4927.125
$ws.Range("M122").Value = -5681.0905
$ws.Range("N122").Value = -9827.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 201828.1
$ws.Range("I105").Value = 127033.75
$ws.Range("K105").Value = 127033.75
$ws.Range("M105").Value = -125286.75
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 66875
$ws.Range("J132").Value = 66875
$ws.Range("L132").Value = 66875
$ws.Range("N132").Value = -76995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 11602.4
$ws.Range("I99").Value = 1443.2
$ws.Range("J99").Value = 21761.6
$ws.Range("K99").Value = 1443.2
$ws.Range("L99").Value = 21761.6
$ws.Range("M99").Value = 54.79999999999995
$ws.Range("N99").Value = -24757.6
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 4448.5713
$ws.Range("I122").Value = 3781.375
$ws.Range("J122").Value = 6583.6
$ws.Range("K122").Value = 11344.125
$ws.Range("L122").Value = 19750.8
$ws.Range("M122").Value = -8894.125
$ws.Range("N122").Value = -24650.8
$ws.Range("H126").Value = 11602.4
$ws.Range("I126").Value = 1443.2
$ws.Range("J126").Value = 21761.6
$ws.Range("K126").Value = 4329.6
$ws.Range("L126").Value = 65284.8
$ws.Range("M126").Value = -1859.6
$ws.Range("N126").Value = -70224.79999999999
$ws.Range("H132").Value = 62505004
$ws.Range("I132").Value = 66672230
$ws.Range("K132").Value = 200016690
$ws.Range("M132").Value = -200014160
$ws.Range("H134").Value = 1425
$ws.Range("I134").Value = 1425
$ws.Range("K134").Value = 4275
$ws.Range("M134").Value = -1740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 115
$ws.Range("J2").Value = 150
$ws.Range("L2").Value = 900
$ws.Range("N2").Value = -1126
$ws.Range("H131").Value = 792.14
$ws.Range("I131").Value = 566.0909
$ws.Range("J131").Value = 820.0787
$ws.Range("K131").Value = 1698.2727
$ws.Range("L131").Value = 2460.2361
$ws.Range("M131").Value = 3341.7273
$ws.Range("N131").Value = -12540.2361
$ws.Range("H136").Value = 2299.2
$ws.Range("J136").Value = 2500
$ws.Range("L136").Value = 7500
$ws.Range("N136").Value = -17700

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2503.0908
$ws.Range("J80").Value = 2491.5
$ws.Range("L80").Value = 2491.5
$ws.Range("N80").Value = -4487.5
$ws.Range("H83").Value = 2503.0908
$ws.Range("J83").Value = 2491.5
$ws.Range("L83").Value = 12457.5
$ws.Range("N83").Value = -22441.5
$ws.Range("H102").Value = 1003398.7
$ws.Range("I102").Value = 4078.4
$ws.Range("J102").Value = 6000000
$ws.Range("K102").Value = 4078.4
$ws.Range("L102").Value = 6000000
$ws.Range("M102").Value = -2456.4
$ws.Range("N102").Value = -6003244

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2195.95
$ws.Range("I7").Value = 1759.9286
$ws.Range("K7").Value = 1759.9286
$ws.Range("M7").Value = -1647.9286
$ws.Range("H40").Value = 79015.766
$ws.Range("J40").Value = 2389.4443
$ws.Range("L40").Value = 2389.4443
$ws.Range("N40").Value = -2661.4443
$ws.Range("H82").Value = 1469.3334
$ws.Range("I82").Value = 1400.6364
$ws.Range("K82").Value = 1400.6364
$ws.Range("M82").Value = -1039.6364
$ws.Range("H85").Value = 1469.3334
$ws.Range("I85").Value = 1400.6364
$ws.Range("K85").Value = 1400.6364
$ws.Range("M85").Value = -152.6364000000001
$ws.Range("H99").Value = 35000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 35000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 35000
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -40990
$ws.Range("H126").Value = 2195.95
$ws.Range("I126").Value = 1759.9286
$ws.Range("K126").Value = 5279.7858
$ws.Range("M126").Value = -2809.7858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5002470
$ws.Range("I62").Value = 25001250
$ws.Range("J62").Value = 2775
$ws.Range("K62").Value = 25001250
$ws.Range("L62").Value = 2775
$ws.Range("M62").Value = -25000626
$ws.Range("N62").Value = -4023
$ws.Range("H65").Value = 5002470
$ws.Range("I65").Value = 25001250
$ws.Range("J65").Value = 2775
$ws.Range("K65").Value = 125006250
$ws.Range("L65").Value = 13875
$ws.Range("M65").Value = -125003130
$ws.Range("N65").Value = -20115
$ws.Range("H122").Value = 1300.5238
$ws.Range("I122").Value = 1107.2
$ws.Range("J122").Value = 1783.8334
$ws.Range("K122").Value = 3321.6
$ws.Range("L122").Value = 5351.5002
$ws.Range("M122").Value = -871.6000000000004
$ws.Range("N122").Value = -10251.5002
$ws.Range("H126").Value = 1224
$ws.Range("I126").Value = 1338.7333
$ws.Range("J126").Value = 879.8
$ws.Range("K126").Value = 4016.199900000001
$ws.Range("L126").Value = 2639.4
$ws.Range("M126").Value = -1546.199900000001
$ws.Range("N126").Value = -7579.4

Write-Host "Applied scheduled market-data refresh updates."
